$d = $word.ActiveDocument

# Commit message: "Als een topic gesloten is, kunnen gebruikers niet meer
# posten" (If a topic is closed, users can no longer post).
#
# In the "Overige taken" table, the rows for "Eigen topic sluiten" and
# "Moderator topic sluiten" each have a remark cell with two paragraphs:
#   "Alle users kunnen dit nog"
#   "Users kunnen nog posten"
# Since a closed topic no longer allows posting, the second paragraph
# ("Users kunnen nog posten") is no longer true and must be removed
# (together with its paragraph mark), leaving just the first paragraph
# in that cell.

$targetText = "Users kunnen nog posten"
$removed = 0

# Walk backwards so deleting a paragraph never invalidates the index of
# a paragraph we still need to inspect.
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text
    # Paragraph text includes the trailing paragraph mark (CR) - strip it
    # for an exact comparison.
    $trimmed = $text.TrimEnd([char]13, [char]7)
    if ($trimmed -eq $targetText) {
        $para.Range.Delete()
        $removed = $removed + 1
    }
}

Write-Output ("Removed " + $removed + " 'Users kunnen nog posten' paragraph(s).")
